# Apply the "bill4" edit to Sheet1:
#  1. Fix E39 (23rd row value) from 470 to 435.
#  2. Fill in the missing S.No. (column A) values for rows 40..347 (A = row-1),
#     matching the existing numbering pattern used in rows 2..39.
#  3. Correct D274 from "isran" to "irsad".
#  4. Move the saved view/selection to the bottom of the list (A325 / F340).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1. Correct the mis-entered amount in row 39.
$ws.Range("E39").Value = 435

# 2. Backfill the running S.No. column for rows 40 through 347.
for ($r = 40; $r -le 347; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# 3. Fix the mistyped name in D274.
$ws.Range("D274").Value = "irsad"

# 4. Update the stored view position/selection (scrolled to the end of the data).
$win = $excel.ActiveWindow
$win.ScrollRow = 325
$win.ScrollColumn = 1
$ws.Range("F340").Select()
